# Updates cryptos list: refreshed Price (col D) and Volume(1h) (col E)
# values for the latest snapshot. Cells in column D whose new text would
# otherwise be auto-parsed as a plain number are entered with a leading
# apostrophe (quote-prefix) so they stay literal text, matching the
# existing inline-string cells in this sheet (e.g. "1.002", "291.26").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "22.441.67"
$ws.Cells.Item(2, 5).Value = "  +0.07%  "
$ws.Cells.Item(3, 4).Value = "1.573.90"
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(6, 4).Value = "'291.14"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "
$ws.Cells.Item(7, 4).Value = "'0.3739"
$ws.Cells.Item(7, 5).Value = "  -0.77%  "
$ws.Cells.Item(8, 4).Value = "'49.93"
$ws.Cells.Item(8, 5).Value = "  +0.08%  "
$ws.Cells.Item(9, 4).Value = "'0.3401"
$ws.Cells.Item(9, 5).Value = "  -0.71%  "
$ws.Cells.Item(10, 4).Value = "'0.07553"
$ws.Cells.Item(10, 5).Value = "  -1.51%  "
$ws.Cells.Item(11, 4).Value = "'1.138"
$ws.Cells.Item(11, 5).Value = "  -2.07%  "
$ws.Cells.Item(12, 4).Value = "'1.003"
$ws.Cells.Item(12, 5).Value = "  +0.01%  "
$ws.Cells.Item(13, 4).Value = "'21.37"
$ws.Cells.Item(13, 5).Value = "  +0.59%  "
$ws.Cells.Item(14, 5).Value = "  -0.25%  "
$ws.Cells.Item(15, 4).Value = "'6.948"
$ws.Cells.Item(15, 5).Value = "  +0.17%  "
$ws.Cells.Item(16, 4).Value = "1.567.97"
$ws.Cells.Item(16, 5).Value = "  -0.35%  "
$ws.Cells.Item(17, 5).Value = "  -0.95%  "
$ws.Cells.Item(18, 4).Value = "'91.08"
$ws.Cells.Item(18, 5).Value = "  +0.82%  "
$ws.Cells.Item(19, 4).Value = "'0.06747"
$ws.Cells.Item(19, 5).Value = "  -0.11%  "
$ws.Cells.Item(20, 5).Value = "  +0.00%  "
$ws.Cells.Item(21, 4).Value = "'6.272"
$ws.Cells.Item(21, 5).Value = "  +0.71%  "
$ws.Cells.Item(22, 5).Value = "  -2.46%  "
$ws.Cells.Item(23, 4).Value = "'12.14"
$ws.Cells.Item(23, 5).Value = "  +0.67%  "
$ws.Cells.Item(24, 4).Value = "22.447.29"
$ws.Cells.Item(24, 5).Value = "  +0.07%  "
$ws.Cells.Item(25, 4).Value = "'2.332"
$ws.Cells.Item(25, 5).Value = "  -3.98%  "
$ws.Cells.Item(26, 4).Value = "'2.590"
$ws.Cells.Item(26, 5).Value = "  -5.33%  "
$ws.Cells.Item(27, 4).Value = "'20.16"
$ws.Cells.Item(27, 5).Value = "  -0.85%  "
$ws.Cells.Item(28, 4).Value = "'148.48"
$ws.Cells.Item(28, 5).Value = "  +1.62%  "
$ws.Cells.Item(29, 4).Value = "'5.010"
$ws.Cells.Item(29, 5).Value = "  -0.36%  "
$ws.Cells.Item(30, 4).Value = "'125.87"
$ws.Cells.Item(30, 5).Value = "  -0.44%  "
$ws.Cells.Item(31, 4).Value = "1.743.59"
$ws.Cells.Item(31, 5).Value = "  -0.27%  "
$ws.Cells.Item(32, 4).Value = "'1.055"
$ws.Cells.Item(32, 5).Value = "  +5.20%  "
$ws.Cells.Item(33, 4).Value = "'6.113"
$ws.Cells.Item(33, 5).Value = "  -1.55%  "
$ws.Cells.Item(34, 4).Value = "'1.982"
$ws.Cells.Item(34, 5).Value = "  -1.77%  "
$ws.Cells.Item(35, 4).Value = "'9.807"
$ws.Cells.Item(35, 5).Value = "  -2.47%  "
$ws.Cells.Item(36, 4).Value = "'0.08416"
$ws.Cells.Item(36, 5).Value = "  -1.92%  "
$ws.Cells.Item(37, 4).Value = "'1.388"
$ws.Cells.Item(37, 5).Value = "  +3.66%  "
$ws.Cells.Item(38, 4).Value = "'0.02465"
$ws.Cells.Item(38, 5).Value = "  -3.29%  "
$ws.Cells.Item(39, 4).Value = "'0.2288"
$ws.Cells.Item(39, 5).Value = "  -1.27%  "
$ws.Cells.Item(40, 4).Value = "'0.06531"
$ws.Cells.Item(40, 5).Value = "  -0.95%  "
$ws.Cells.Item(41, 4).Value = "'5.459"
$ws.Cells.Item(41, 5).Value = "  -0.31%  "
$ws.Cells.Item(42, 4).Value = "'11.27"
$ws.Cells.Item(42, 5).Value = "  -2.55%  "
$ws.Cells.Item(43, 4).Value = "'0.6250"
$ws.Cells.Item(43, 5).Value = "  -3.18%  "
$ws.Cells.Item(44, 5).Value = "  -0.04%  "
$ws.Cells.Item(45, 4).Value = "'14.01"
$ws.Cells.Item(45, 5).Value = "  -1.36%  "
$ws.Cells.Item(46, 4).Value = "'3.813"
$ws.Cells.Item(46, 5).Value = "  +0.36%  "
$ws.Cells.Item(47, 4).Value = "'0.5816"
$ws.Cells.Item(47, 5).Value = "  -3.30%  "
$ws.Cells.Item(48, 4).Value = "'2.084"
$ws.Cells.Item(48, 5).Value = "  -0.13%  "
$ws.Cells.Item(49, 4).Value = "'129.31"
$ws.Cells.Item(49, 5).Value = "  +2.95%  "
$ws.Cells.Item(50, 4).Value = "'1.223"
$ws.Cells.Item(50, 5).Value = "  -5.76%  "
$ws.Cells.Item(51, 4).Value = "'0.07326"
